$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '23.177.52'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.602.12'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '303.25'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3783'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '52.37'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.85%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3615'
$ws.Range('D9').Style = 'Normal'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.265'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.87%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.002'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.02%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.08130'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.11%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.64'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.581'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.389'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001247'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.601.68'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.08%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '93.84'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06882'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.97%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.552'
$ws.Range('D21').Style = 'Normal'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.94'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.16%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '23.185.80'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.45%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.386'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +2.10%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.983'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.22%  '
$ws.Range('E27').Value = '  +0.85%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '148.89'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.49%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.256'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.30%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '133.75'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.64%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.375'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.72%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '6.827'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.778.55'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9724'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +2.14%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.07521'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.45%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.32'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02721'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.18%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.2509'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.28%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.119'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.08815'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.84%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.363'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.59%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.7087'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '12.51'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.50%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '15.55'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6529'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.53%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.309'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.53%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.010'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '132.14'
$ws.Range('D48').Style = 'Normal'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.07954'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.28%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.203'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.01%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.219'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.34%  '
